$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.106.28'
$ws.Range('E2').Value = '  +2.87%  '
$ws.Range('D3').Value = '2.957.95'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '595.21'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '149.02'
$ws.Range('E6').Value = '  +2.75%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '2.955.44'
$ws.Range('E8').Value = '  +1.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.508'
$ws.Range('E9').Value = '  +1.38%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.15'
$ws.Range('E10').Value = '  +3.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.151'
$ws.Range('E11').Value = '  +6.73%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.442'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000236'
$ws.Range('E13').Value = '  +5.21%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.87'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').Value = '3.446.34'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').Value = '63.032.28'
$ws.Range('E17').Value = '  +2.87%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.71'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').Value = '2.956.87'
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '442.50'
$ws.Range('E20').Value = '  +2.40%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.53'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.670'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.02'
$ws.Range('E23').Value = '  -0.83%  '
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.22'
$ws.Range('E24').Value = '  +2.51%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '81.04'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.15'
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.80'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.32'
$ws.Range('E29').Value = '  +6.14%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.21'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0000102'
$ws.Range('E32').Value = '  +15.79%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.51'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('E34').Value = '  -0.69%  '
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.993'
$ws.Range('E36').Value = '  -1.66%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.62'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.09'
$ws.Range('E38').Value = '  +3.72%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.05'
$ws.Range('E39').Value = '  +2.49%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '49.73'
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '8.52'
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.118'
$ws.Range('E42').Value = '  -4.44%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.281'
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '39.16'
$ws.Range('E44').Value = '  -7.75%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.699.22'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '135.32'
$ws.Range('E46').Value = '  +1.31%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0338'
$ws.Range('E47').Value = '  -1.85%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '360.63'
$ws.Range('E48').Value = '  -1.62%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.105'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '22.84'
$ws.Range('E51').Value = '  -2.93%  '
